$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.389.07"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").Value = "1.852.80"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.74%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.11"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6200"
$ws.Range("E6").Value = "  -1.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.015"
$ws.Range("E7").Value = "  +1.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07478"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2964"
$ws.Range("E9").Value = "  +1.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.12"
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  +0.37%  "

# Row 12
$ws.Range("D12").Value = "1.830.90"
$ws.Range("E12").Value = "  -0.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.029"
$ws.Range("E13").Value = "  +0.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6768"
$ws.Range("E14").Value = "  +1.19%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.36"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009085"
$ws.Range("E16").Value = "  -3.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.916"
$ws.Range("E17").Value = "  -2.64%  "

# Row 18
$ws.Range("D18").Value = "29.337.77"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19
$ws.Range("D19").Value = "2.083.63"
$ws.Range("E19").Value = "  +0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.69"
$ws.Range("E20").Value = "  +6.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.73"
$ws.Range("E21").Value = "  +0.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.017"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.204"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.014"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.74"
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1440"
$ws.Range("E26").Value = "  +2.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.557"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.97"
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  +0.81%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05637"
$ws.Range("E30").Value = "  +2.88%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.173"
$ws.Range("E31").Value = "  +0.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.129"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.225"
$ws.Range("E33").Value = "  +1.64%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.860"
$ws.Range("E34").Value = "  +0.20%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7490"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.148"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.671"
$ws.Range("E37").Value = "  +2.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.834"
$ws.Range("E38").Value = "  +2.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01790"
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("D40").Value = "1.219.76"
$ws.Range("E40").Value = "  -0.72%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.516"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9078"
$ws.Range("E42").Value = "  +1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.016"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.80"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("D45").Value = "1.984.13"
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.56"
$ws.Range("E46").Value = "  -0.14%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  -1.43%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5143"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4071"
$ws.Range("E49").Value = "  +0.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.174"
$ws.Range("E50").Value = "  +1.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05850"
$ws.Range("E51").Value = "  +0.82%  "
